# Edit: insert a new week's worth of Coliflor pricing data (2 rows) at the
# top of the Femacal de La Calera data block (rows 275-276), pushing all
# subsequent rows down by two. The two newly inserted rows contain a new
# observation for date 44489 (Primera / Segunda quality).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 275; this shifts the existing rows 275-395
# down to 277-397 (and all their formatting/values) automatically.
$ws.Rows("275:276").Insert()

# --- New row 275 (Primera) ---
$ws.Cells.Item(275, 1).Value = 3
$ws.Cells.Item(275, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(275, 3).Value = "Coquimbo"
$ws.Cells.Item(275, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(275, 4).Value = 44489
$ws.Cells.Item(275, 5).Value = 5
$ws.Cells.Item(275, 6).Value = 100112008
$ws.Cells.Item(275, 7).Value = "Coliflor"
$ws.Cells.Item(275, 8).Value = "Sin especificar"
$ws.Cells.Item(275, 9).Value = "Primera"
$ws.Cells.Item(275, 10).Value = 1850
$ws.Cells.Item(275, 11).Value = 600
$ws.Cells.Item(275, 12).Value = 650
$ws.Cells.Item(275, 13).Value = 626
$ws.Cells.Item(275, 14).Value = "`$/unidad"
$ws.Cells.Item(275, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(275, 16).Value = 626
$ws.Cells.Item(275, 17).Value = 1
$ws.Cells.Item(275, 18).Value = "Hortaliza"

# --- New row 276 (Segunda) ---
$ws.Cells.Item(276, 1).Value = 3
$ws.Cells.Item(276, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(276, 3).Value = "Coquimbo"
$ws.Cells.Item(276, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(276, 4).Value = 44489
$ws.Cells.Item(276, 5).Value = 5
$ws.Cells.Item(276, 6).Value = 100112008
$ws.Cells.Item(276, 7).Value = "Coliflor"
$ws.Cells.Item(276, 8).Value = "Sin especificar"
$ws.Cells.Item(276, 9).Value = "Segunda"
$ws.Cells.Item(276, 10).Value = 850
$ws.Cells.Item(276, 11).Value = 500
$ws.Cells.Item(276, 12).Value = 500
$ws.Cells.Item(276, 13).Value = 500
$ws.Cells.Item(276, 14).Value = "`$/unidad"
$ws.Cells.Item(276, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(276, 16).Value = 500
$ws.Cells.Item(276, 17).Value = 1
$ws.Cells.Item(276, 18).Value = "Hortaliza"
